# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1;  3 = 0;  4 = 2;  5 = 1;  6 = 0;  7 = 0;  8 = 1;  9 = 1;  10 = 1;
    11 = 2; 12 = 1; 13 = 0; 14 = 2; 15 = 1; 16 = 1; 17 = 1; 18 = 2; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 0; 25 = 0; 26 = 2; 27 = 1; 28 = 0;
    29 = 0; 30 = 0; 31 = 1; 32 = 2; 33 = 3; 34 = 2; 35 = 1; 36 = 2; 37 = 2;
    38 = 2; 39 = 0; 40 = 1; 41 = 0; 42 = 2; 43 = 3; 44 = 1; 45 = 0; 46 = 0;
    47 = 1; 48 = 0; 49 = 1; 50 = 2; 51 = 2; 52 = 0; 53 = 1; 54 = 0; 55 = 1;
    56 = 0; 57 = 1; 58 = 1; 59 = 2; 60 = 2; 61 = 0; 62 = 0; 63 = 0; 64 = 2;
    65 = 0; 66 = 3; 67 = 0; 68 = 0; 69 = 1; 70 = 0; 71 = 0; 72 = 2; 73 = 2;
    74 = 3; 75 = 1; 76 = 1; 77 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
